# Cleaned up BasicResultsLoad test
# - Adds an "Invalid (age)" regression row to AdvancedResultsLoad
# - Populates BasicResultsLoad with its full set of Params/SearchType/
#   CtsFieldsUsed/CtsCancerInfo/CtsLocation test rows

$wb = $excel.ActiveWorkbook

$wsAdvanced = $wb.Worksheets.Item("AdvancedResultsLoad")
$wsBasic    = $wb.Worksheets.Item("BasicResultsLoad")

# ---------------------------------------------------------------------
# BasicResultsLoad: fill in the full test matrix (rows 2-15)
# ---------------------------------------------------------------------
$basicRows = @(
    @("?q=&t=&a=&z=&rl=1",                     "All Trials",              "none",      "none|none",                  "all"),
    @("?q=&t=&a=&z=&rl=1",                     "All Trials (no params)",  "none",      "none|none",                  "all"),
    @("?t=C9087&a=&z=&rl=1",                   "Cancer Type",             "t",         "typecondition|c9087|none",   "all"),
    @("?q=medulla&t=&a=&z=&rl=1",               "Keyword",                 "q",         "keyword|medulla|none",       "all"),
    @("?q=&t=&a=80&z=&rl=1",                    "Age",                     "a",         "none|80",                    "all"),
    @("?q=&t=&a=&z=20772&rl=1",                 "Zip",                     "loc:z",     "none|none",                  "zip|20772|none"),
    @("?t=C3869&a=85&z=&rl=1",                  "Cancer Type/Age",         "t:a",       "typecondition|c3869|85",     "all"),
    @("?q=Interstitial+&t=&a=78&z=&rl=1",        "Keyword/Age",             "a:q",       "keyword|interstitial|78",    "all"),
    @("?t=C7853&a=&z=29401&rl=1",                "Cancer Type/Zip",         "t:loc:z",   "typecondition|c7853|none",   "zip|29401|none"),
    @("?q=androgen&t=&a=&z=25063&rl=1",          "Keyword/Zip",             "q:loc:z",   "keyword|androgen|none",      "zip|25063|none"),
    @("?t=C3242&a=101&z=96795&rl=1",             "Cancer Type/Age/Zip",     "t:a:loc:z", "typecondition|c3242|101",    "zip|96795|none"),
    @("?q=plasma&t=&a=70&z=60044&rl=1",          "Keyword/Age/Zip",         "a:q:loc:z", "keyword|plasma|70",          "zip|60044|none"),
    @("?q=&t=&a=&z=abcd&rl=1",                   "Invalid (zip)",           "loc",       "none|none",                  "zip|none|none"),
    @("?q=&t=&a=500&z=&rl=1",                    "Invalid (age)",           "none",      "none|none",                  "all")
)

$r = 2
foreach ($row in $basicRows) {
    $wsBasic.Cells.Item($r, 1).Value = $row[0]
    $wsBasic.Cells.Item($r, 2).Value = $row[1]
    $wsBasic.Cells.Item($r, 3).Value = $row[2]
    $wsBasic.Cells.Item($r, 4).Value = $row[3]
    $wsBasic.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$wsBasic.Columns.Item(1).AutoFit()
$wsBasic.Columns.Item(2).AutoFit()
$wsBasic.Columns.Item(4).AutoFit()
$wsBasic.Columns.Item(5).AutoFit()
$wsBasic.Range("A16").Select()

# ---------------------------------------------------------------------
# AdvancedResultsLoad: append the matching "Invalid (age)" row
# ---------------------------------------------------------------------
$wsAdvanced.Cells.Item(22, 1).Value = "?t=&a=age&q=&loc=0&tt=&tp=&tid=&in=&lo=&rl=2"
$wsAdvanced.Cells.Item(22, 2).Value = "Invalid (age)"
$wsAdvanced.Cells.Item(22, 3).Value = "none"
$wsAdvanced.Cells.Item(22, 4).Value = "all|all|all|all|none|none"
$wsAdvanced.Cells.Item(22, 5).Value = "all"
$wsAdvanced.Cells.Item(22, 6).Value = "all|none|none"
$wsAdvanced.Cells.Item(22, 7).Value = "all|none|none|none"

$wsAdvanced.Columns.Item(4).AutoFit()
$wsAdvanced.Columns.Item(5).AutoFit()
$wsAdvanced.Columns.Item(6).AutoFit()

# Restore AdvancedResultsLoad as the active tab/selection (A23, one past
# the newly-added last row), matching the original authoring session.
$wsAdvanced.Activate()
$wsAdvanced.Range("A23").Select()
